# feat: add 2022-Q3 data
#
# - Inserts a new "2022-Q3" worksheet (fund-holdings detail) right before
#   the existing "2022-Q2" sheet, cloned from "2022-Q2" so it inherits the
#   same header/index-column formatting, then repopulated with the Q3 data.
# - Updates the "总计" roll-up sheet: the old rows shift down one slot and
#   a new top row for 2022-Q3 is inserted.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $value) {
    # Force a literal text value even when it looks like a number
    # (e.g. "202001", "19.56"), matching the source file's inlineStr
    # cells, without leaving a stray NumberFormat-driven style behind.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# ---------------------------------------------------------------------
# 1. "总计" roll-up sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Clone the existing index-column style (bold + border, s="2") down into
# the brand-new row 6.
$summary.Range("A5").Copy($summary.Range("A6"))

$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2020-Q4"
$summary.Range("C6").Value = 2
$summary.Range("D6").Value = 0.07000000000000001

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2021-Q4"
$summary.Range("C5").Value = 6
$summary.Range("D5").Value = 0.25

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 2
$summary.Range("D4").Value = 0.03

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 4
$summary.Range("D3").Value = 0.05

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 12
$summary.Range("D2").Value = 1.63

# ---------------------------------------------------------------------
# 2. New "2022-Q3" detail sheet
# ---------------------------------------------------------------------
# A whole-sheet copy (rather than adding a blank sheet + pasting ranges
# into it) is used so the new sheet reliably inherits "2022-Q2"'s cell
# formatting (bold/bordered header row + index column).
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# Extend the bold/bordered index-column style from row 5 down through the
# extra rows this quarter needs (13 total vs. 5 in 2022-Q2).
for ($i = 6; $i -le 13; $i++) {
    $q3.Range("A5").Copy($q3.Range("A$i"))
}

$rows = @(
    @(0,  "202001", "南方稳健成长混合",            "19.56", "76.28", "1.88", "0.3677", 7),
    @(1,  "160325", "华夏创业板两年定期开放混合",  "8.96",  "94.23", "4.02", "0.3602", 6),
    @(2,  "202002", "南方稳健成长贰号混合",        "14.12", "76.88", "1.88", "0.2655", 7),
    @(3,  "011216", "南方优质企业混合A",           "7.95",  "92.17", "1.94", "0.1542", 10),
    @(4,  "163302", "大摩资源优选混合（LOF）",     "5.08",  "79.56", "2.92", "0.1483", 10),
    @(5,  "160143", "南方创业板2年定期开放混合",   "3.38",  "83.97", "4.12", "0.1393", 3),
    @(6,  "002160", "南方转型驱动灵活配置混合",    "3.14",  "93.43", "2.81", "0.0882", 5),
    @(7,  "009847", "圆信永丰研究精选混合A",       "1.14",  "89.51", "3.38", "0.0385", 8),
    @(8,  "006969", "圆信永丰高端制造混合",        "0.91",  "87.79", "3.29", "0.0299", 7),
    @(9,  "011217", "南方优质企业混合C",           "1.18",  "92.17", "1.94", "0.0229", 10),
    @(10, "009848", "圆信永丰研究精选混合C",       "0.44",  "89.51", "3.38", "0.0149", 8),
    @(11, "009054", "圆信永丰沣泰混合",            "0.23",  "26.81", "1.12", "0.0026", 7)
)

$r = 2
foreach ($row in $rows) {
    $q3.Range("A$r").Value = $row[0]
    Set-TextValue $q3.Range("B$r") $row[1]
    Set-TextValue $q3.Range("C$r") $row[2]
    Set-TextValue $q3.Range("D$r") $row[3]
    Set-TextValue $q3.Range("E$r") $row[4]
    Set-TextValue $q3.Range("F$r") $row[5]
    Set-TextValue $q3.Range("G$r") $row[6]
    $q3.Range("H$r").Value = $row[7]
    $r = $r + 1
}
